$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 13 new rows of generated question/answer data (rows 23-35)
$ws.Range("B23").Value = 'Why'
$ws.Range("C23").Value = 'Why might the most anterior face patch (AM) show a weaker representation of non-face stimuli than other face patches?'
$ws.Range("D23").Value = 'Faces are a behaviorally important class of visual stimuli for primates. Recent work in macaque monkeys has identified six discrete face areas where most neurons have higher firing rates to images of faces compared with other objects (Tsao et al., 2006). While neurons in these areas appear to have different tuning (Freiwald and Tsao, 2010; Issa and DiCarlo, 2012), exactly what types of information and, consequently, which visual behaviors neural populations within each face area can support, is unknown. Here we use population decoding to better characterize three of these face patches (ML/MF, AL, and AM). We show that neural activity in all patches contains information that discriminates between the broad categories of face and nonface objects, individual faces, and nonface stimuli. Information is present in both high and lower firing rate regimes. However, there were significant differences between the patches, with the most anterior patch showing relatively weaker representation of nonface stimuli. Additionally, we find that pose-invariant face identity information increases as one moves to more anterior patches, while information about the orientation of the head decreases. Finally, we show that all the information we can extract from the population is present in patterns of activity across neurons, and there is relatively little information in the total activity of the population. These findings give new insight into the representations constructed by the face patch system and how they are successively transformed.'

$ws.Range("B24").Value = 'How'
$ws.Range("C24").Value = 'How do researchers use population decoding to characterize the neural activity within face patches?'
$ws.Range("D24").Value = 'Faces are a behaviorally important class of visual stimuli for primates. Recent work in macaque monkeys has identified six discrete face areas where most neurons have higher firing rates to images of faces compared with other objects (Tsao et al., 2006). While neurons in these areas appear to have different tuning (Freiwald and Tsao, 2010; Issa and DiCarlo, 2012), exactly what types of information and, consequently, which visual behaviors neural populations within each face area can support, is unknown. Here we use population decoding to better characterize three of these face patches (ML/MF, AL, and AM). We show that neural activity in all patches contains information that discriminates between the broad categories of face and nonface objects, individual faces, and nonface stimuli. Information is present in both high and lower firing rate regimes. However, there were significant differences between the patches, with the most anterior patch showing relatively weaker representation of nonface stimuli. Additionally, we find that pose-invariant face identity information increases as one moves to more anterior patches, while information about the orientation of the head decreases. Finally, we show that all the information we can extract from the population is present in patterns of activity across neurons, and there is relatively little information in the total activity of the population. These findings give new insight into the representations constructed by the face patch system and how they are successively transformed.'

$ws.Range("B25").Value = 'Hypothetical'
$ws.Range("C25").Value = 'If a patient with schizophrenia presents with low IQ and severe negative symptoms, would they be more likely to demonstrate insufficient effort on neuropsychological testing, even if positive symptoms are less pronounced?'
$ws.Range("D25").Value = 'There is some evidence that insufficient effort may be common in schizophrenia, posing significant threats to the validity of neuropsychological test results. Low effort may account for a significant proportion of variance in neuropsychological test scores and the generalized cognitive deficit that characterizes the disorder. The current study evaluated clinical predictors of insufficient effort in schizophrenia using an embedded effort measure, the Repeatable Battery for the Assessment of Neuropsychological Status (RBANS) Effort Index (EI). Participants were 330 patients meeting DSM-IV-TR criteria for schizophrenia, schizoaffective disorder, or another psychotic disorder who received a battery of neuropsychological tests, including: Wechsler Test of Adult Reading (WTAR), Wechsler Abbreviated Scale of Intelligence (WASI), and RBANS. Clinical assessments designed to measure functional outcome and symptoms were also obtained. Results indicated that 9.4% of patients failed the EI. Patients who failed had lower full-scale, verbal, and performance IQ, as well as poorer performance on RBANS domains not included in the EI (immediate memory, language, and visuospatial/construction). Patients who failed the EI also displayed poorer community-based vocational outcome, greater likelihood of having "deficit schizophrenia" (i.e., primary and enduring negative symptoms), and increased severity of positive symptoms. Regression analyses revealed that insufficient effort was most significantly predicted by a combination of low IQ, negative symptoms, and positive symptoms. Findings suggest that although insufficient effort may be relatively uncommon in schizophrenia, it is associated with important clinical outcomes. The RBANS EI may be a useful tool in evaluating insufficient effort in schizophrenia.'

$ws.Range("B26").Value = 'Hypothetical'
$ws.Range("C26").Value = 'If the prevalence of "deficit schizophrenia" increases within a clinical population, might that lead to a corresponding increase in insufficient effort during neuropsychological evaluations?'
$ws.Range("D26").Value = 'There is some evidence that insufficient effort may be common in schizophrenia, posing significant threats to the validity of neuropsychological test results. Low effort may account for a significant proportion of variance in neuropsychological test scores and the generalized cognitive deficit that characterizes the disorder. The current study evaluated clinical predictors of insufficient effort in schizophrenia using an embedded effort measure, the Repeatable Battery for the Assessment of Neuropsychological Status (RBANS) Effort Index (EI). Participants were 330 patients meeting DSM-IV-TR criteria for schizophrenia, schizoaffective disorder, or another psychotic disorder who received a battery of neuropsychological tests, including: Wechsler Test of Adult Reading (WTAR), Wechsler Abbreviated Scale of Intelligence (WASI), and RBANS. Clinical assessments designed to measure functional outcome and symptoms were also obtained. Results indicated that 9.4% of patients failed the EI. Patients who failed had lower full-scale, verbal, and performance IQ, as well as poorer performance on RBANS domains not included in the EI (immediate memory, language, and visuospatial/construction). Patients who failed the EI also displayed poorer community-based vocational outcome, greater likelihood of having "deficit schizophrenia" (i.e., primary and enduring negative symptoms), and increased severity of positive symptoms. Regression analyses revealed that insufficient effort was most significantly predicted by a combination of low IQ, negative symptoms, and positive symptoms. Findings suggest that although insufficient effort may be relatively uncommon in schizophrenia, it is associated with important clinical outcomes. The RBANS EI may be a useful tool in evaluating insufficient effort in schizophrenia.'

$ws.Range("B27").Value = 'Open Ended Question'
$ws.Range("D27").Value = 'Using a "theory of mind" allows us to explain and predict others'' behavior in terms of their mental states, yet individual differences in the accuracy of mental state inferences are not well understood. We hypothesized that the accuracy of mental state inferences can be explained by the ability to characterize the mind giving rise to the mental state. Under this proposal, individuals differentiate between minds by representing them in "Mind-space"-a multidimensional space where dimensions reflect any characteristic of minds that allows them to be individuated. Individual differences in the representation of minds and the accuracy of mental state inferences are explained by one''s model of how minds can vary (Mind-space) and ability to locate an individual mind within this space. We measured the accuracy of participants'' model of the covariance between dimensions in Mind-space that represent personality traits, and we found this was associated with the accuracy of mental state inference (Experiment 1). Mind-space accuracy also predicted the ability to locate others within Mind-space on dimensions of personality and intelligence (Experiment 2). Direct evidence for the representation of minds in mental state inference was obtained by showing that the location of others in Mind-space affects the probability of particular mental states being ascribed to them (Experiment 3). This latter effect extended to mental states dependent upon representation of trait covariation (Experiment 4). Results support the claim that mental state inference varies according to location in Mind-space, and therefore that adopting the Mind-space framework can explain some of the individual differences in theory of mind. (PsycInfo Database Record (c) 2020 APA, all rights reserved)'
$ws.Range("C27").Value = 'What factors, besides personality traits and intelligence, might be used as dimensions within an individual''s "Mind-space"?'

$ws.Range("B28").Value = 'Yes / No question'
$ws.Range("C28").Value = 'Can an individual''s "Mind-space" model change or evolve over time?'
$ws.Range("D28").Value = 'Using a "theory of mind" allows us to explain and predict others'' behavior in terms of their mental states, yet individual differences in the accuracy of mental state inferences are not well understood. We hypothesized that the accuracy of mental state inferences can be explained by the ability to characterize the mind giving rise to the mental state. Under this proposal, individuals differentiate between minds by representing them in "Mind-space"-a multidimensional space where dimensions reflect any characteristic of minds that allows them to be individuated. Individual differences in the representation of minds and the accuracy of mental state inferences are explained by one''s model of how minds can vary (Mind-space) and ability to locate an individual mind within this space. We measured the accuracy of participants'' model of the covariance between dimensions in Mind-space that represent personality traits, and we found this was associated with the accuracy of mental state inference (Experiment 1). Mind-space accuracy also predicted the ability to locate others within Mind-space on dimensions of personality and intelligence (Experiment 2). Direct evidence for the representation of minds in mental state inference was obtained by showing that the location of others in Mind-space affects the probability of particular mental states being ascribed to them (Experiment 3). This latter effect extended to mental states dependent upon representation of trait covariation (Experiment 4). Results support the claim that mental state inference varies according to location in Mind-space, and therefore that adopting the Mind-space framework can explain some of the individual differences in theory of mind. (PsycInfo Database Record (c) 2020 APA, all rights reserved)'

$ws.Range("B29").Value = 'Multiple choice question'
$ws.Range("C29").Value = 'According to the study, which of the following is the primary factor influencing the accuracy of mental state inferences?
a) The number of dimensions used in an individual''s "Mind-space" model
b) The accuracy of an individual''s model of how minds covary
c) The ability to locate a specific mind within one''s "Mind-space"
d) Pre-existing biases about other people''s mental states'
$ws.Range("D29").Value = 'Using a "theory of mind" allows us to explain and predict others'' behavior in terms of their mental states, yet individual differences in the accuracy of mental state inferences are not well understood. We hypothesized that the accuracy of mental state inferences can be explained by the ability to characterize the mind giving rise to the mental state. Under this proposal, individuals differentiate between minds by representing them in "Mind-space"-a multidimensional space where dimensions reflect any characteristic of minds that allows them to be individuated. Individual differences in the representation of minds and the accuracy of mental state inferences are explained by one''s model of how minds can vary (Mind-space) and ability to locate an individual mind within this space. We measured the accuracy of participants'' model of the covariance between dimensions in Mind-space that represent personality traits, and we found this was associated with the accuracy of mental state inference (Experiment 1). Mind-space accuracy also predicted the ability to locate others within Mind-space on dimensions of personality and intelligence (Experiment 2). Direct evidence for the representation of minds in mental state inference was obtained by showing that the location of others in Mind-space affects the probability of particular mental states being ascribed to them (Experiment 3). This latter effect extended to mental states dependent upon representation of trait covariation (Experiment 4). Results support the claim that mental state inference varies according to location in Mind-space, and therefore that adopting the Mind-space framework can explain some of the individual differences in theory of mind. (PsycInfo Database Record (c) 2020 APA, all rights reserved)'

$ws.Range("B30").Value = 'Open-ended question'
$ws.Range("D30").Value = 'Intelligent life has emerged late in Earth''s habitable lifetime, and required a preceding series of key evolutionary transitions. A simple model (the Carter model) explains the late arrival of intelligent life by positing these evolutionary transitions were exceptionally unlikely ''critical steps''. An alternative model (the neocatastrophism hypothesis) proposes that intelligent life was delayed by frequent catastrophes that served to set back evolutionary innovation. Here, we generalize the Carter model and explore this hypothesis by including catastrophes that can ''undo'' an evolutionary transition. Introducing catastrophes or evolutionary dead ends can create situations in which critical steps occur rapidly or in clusters, suggesting that past estimates of the number of critical steps could be underestimated. If catastrophes affect complex life more than simple life, the critical steps will also exhibit a pattern of acceleration towards the present, suggesting that the increase in biological complexity over the past 500 Myr could reflect previously overlooked evolutionary transitions. Furthermore, our results have implications for understanding the different explanations (critical steps versus neo-catastrophes) for the evolution of intelligent life and the so-called Fermi paradox-the observation that intelligent life appears rare in the observable Universe.'
$ws.Range("C30").Value = 'How does the revised Carter model, which includes catastrophes, change our understanding of the "critical steps" involved in the evolution of intelligent life?'

$ws.Range("B31").Value = 'Yes / No question'
$ws.Range("C31").Value = 'Could the apparent increase in biological complexity over the past 500 million years be partially explained by evolutionary transitions that researchers previously hadn''t considered?'
$ws.Range("D31").Value = 'Intelligent life has emerged late in Earth''s habitable lifetime, and required a preceding series of key evolutionary transitions. A simple model (the Carter model) explains the late arrival of intelligent life by positing these evolutionary transitions were exceptionally unlikely ''critical steps''. An alternative model (the neocatastrophism hypothesis) proposes that intelligent life was delayed by frequent catastrophes that served to set back evolutionary innovation. Here, we generalize the Carter model and explore this hypothesis by including catastrophes that can ''undo'' an evolutionary transition. Introducing catastrophes or evolutionary dead ends can create situations in which critical steps occur rapidly or in clusters, suggesting that past estimates of the number of critical steps could be underestimated. If catastrophes affect complex life more than simple life, the critical steps will also exhibit a pattern of acceleration towards the present, suggesting that the increase in biological complexity over the past 500 Myr could reflect previously overlooked evolutionary transitions. Furthermore, our results have implications for understanding the different explanations (critical steps versus neo-catastrophes) for the evolution of intelligent life and the so-called Fermi paradox-the observation that intelligent life appears rare in the observable Universe.'

$ws.Range("B32").Value = 'Comparative'
$ws.Range("C32").Value = 'Which hypothesis (critical steps or neo-catastrophes) provides a more plausible explanation for the Fermi paradox?'
$ws.Range("D32").Value = 'Intelligent life has emerged late in Earth''s habitable lifetime, and required a preceding series of key evolutionary transitions. A simple model (the Carter model) explains the late arrival of intelligent life by positing these evolutionary transitions were exceptionally unlikely ''critical steps''. An alternative model (the neocatastrophism hypothesis) proposes that intelligent life was delayed by frequent catastrophes that served to set back evolutionary innovation. Here, we generalize the Carter model and explore this hypothesis by including catastrophes that can ''undo'' an evolutionary transition. Introducing catastrophes or evolutionary dead ends can create situations in which critical steps occur rapidly or in clusters, suggesting that past estimates of the number of critical steps could be underestimated. If catastrophes affect complex life more than simple life, the critical steps will also exhibit a pattern of acceleration towards the present, suggesting that the increase in biological complexity over the past 500 Myr could reflect previously overlooked evolutionary transitions. Furthermore, our results have implications for understanding the different explanations (critical steps versus neo-catastrophes) for the evolution of intelligent life and the so-called Fermi paradox-the observation that intelligent life appears rare in the observable Universe.'

$ws.Range("B33").Value = 'Open-ended question'
$ws.Range("D33").Value = 'Delay discounting (DD) refers to how rapidly an individual devalues goods based on delays to receipt. DD usually is considered a trait variable but can be state dependent, yet few studies have assessed commodity valuation at short, naturalistically relevant time intervals that might enable state-dependent analysis. This study aimed to determine whether drug-use impulsivity and intelligence influence heroin DD at short (ecologically relevant) delays during two pharmacological states (heroin satiation and withdrawal). Out-of-treatment, intensive heroin users (n = 170; 53.5% African American; 66.7% male) provided complete DD data during imagined heroin satiation and withdrawal. Delays were 3, 6, 12, 24, 48, 72, and 96 hours; maximum delayed heroin amount was thirty $10 bags. Indifference points were used to calculate area under the curve (AUC). We also assessed drug-use impulsivity (subscales from the Impulsive Relapse Questionnaire [IRQ]) and estimated intelligence (Shipley IQ) as predictors of DD. Heroin discounting was greater (smaller AUC) during withdrawal than satiation. In regression analyses, lower intelligence and IRQ Capacity for Delay as well as higher IRQ Speed (to return to drug use) predicted greater heroin discounting in the satiation condition. Lower intelligence and higher IRQ Speed predicted greater discounting in the withdrawal condition. Sex, race, substance use variables, and other IRQ subscales were not significantly related to the withdrawal or satiation DD behavior. In summary, heroin discounting was temporally rapid, pharmacologically state dependent, and predicted by drug-use impulsivity and estimated intelligence. These findings highlight a novel and sensitive measure of acute DD that is easy to administer.'
$ws.Range("C33").Value = 'Besides drug-use impulsivity and intelligence, what other factors might contribute to differences in heroin delay discounting (DD) during states of satiation and withdrawal?'

$ws.Range("B34").Value = 'Yes / No question'
$ws.Range("C34").Value = 'Do the factors predicting heroin DD in a state of satiation differ completely from those predicting DD in a state of withdrawal?'
$ws.Range("D34").Value = 'Delay discounting (DD) refers to how rapidly an individual devalues goods based on delays to receipt. DD usually is considered a trait variable but can be state dependent, yet few studies have assessed commodity valuation at short, naturalistically relevant time intervals that might enable state-dependent analysis. This study aimed to determine whether drug-use impulsivity and intelligence influence heroin DD at short (ecologically relevant) delays during two pharmacological states (heroin satiation and withdrawal). Out-of-treatment, intensive heroin users (n = 170; 53.5% African American; 66.7% male) provided complete DD data during imagined heroin satiation and withdrawal. Delays were 3, 6, 12, 24, 48, 72, and 96 hours; maximum delayed heroin amount was thirty $10 bags. Indifference points were used to calculate area under the curve (AUC). We also assessed drug-use impulsivity (subscales from the Impulsive Relapse Questionnaire [IRQ]) and estimated intelligence (Shipley IQ) as predictors of DD. Heroin discounting was greater (smaller AUC) during withdrawal than satiation. In regression analyses, lower intelligence and IRQ Capacity for Delay as well as higher IRQ Speed (to return to drug use) predicted greater heroin discounting in the satiation condition. Lower intelligence and higher IRQ Speed predicted greater discounting in the withdrawal condition. Sex, race, substance use variables, and other IRQ subscales were not significantly related to the withdrawal or satiation DD behavior. In summary, heroin discounting was temporally rapid, pharmacologically state dependent, and predicted by drug-use impulsivity and estimated intelligence. These findings highlight a novel and sensitive measure of acute DD that is easy to administer.'

$ws.Range("B35").Value = 'Multiple choice question'
$ws.Range("C35").Value = 'Which of the following best summarizes the primary finding of the study?
a) Heroin delay discounting is not influenced by pharmacological state.
b) Delay discounting is a trait that remains constant regardless of circumstances.
c) Heroin delay discounting is influenced by both pharmacological state and individual characteristics.
d) Intelligence is the single greatest predictor of heroin delay discounting.'
$ws.Range("D35").Value = 'Delay discounting (DD) refers to how rapidly an individual devalues goods based on delays to receipt. DD usually is considered a trait variable but can be state dependent, yet few studies have assessed commodity valuation at short, naturalistically relevant time intervals that might enable state-dependent analysis. This study aimed to determine whether drug-use impulsivity and intelligence influence heroin DD at short (ecologically relevant) delays during two pharmacological states (heroin satiation and withdrawal). Out-of-treatment, intensive heroin users (n = 170; 53.5% African American; 66.7% male) provided complete DD data during imagined heroin satiation and withdrawal. Delays were 3, 6, 12, 24, 48, 72, and 96 hours; maximum delayed heroin amount was thirty $10 bags. Indifference points were used to calculate area under the curve (AUC). We also assessed drug-use impulsivity (subscales from the Impulsive Relapse Questionnaire [IRQ]) and estimated intelligence (Shipley IQ) as predictors of DD. Heroin discounting was greater (smaller AUC) during withdrawal than satiation. In regression analyses, lower intelligence and IRQ Capacity for Delay as well as higher IRQ Speed (to return to drug use) predicted greater heroin discounting in the satiation condition. Lower intelligence and higher IRQ Speed predicted greater discounting in the withdrawal condition. Sex, race, substance use variables, and other IRQ subscales were not significantly related to the withdrawal or satiation DD behavior. In summary, heroin discounting was temporally rapid, pharmacologically state dependent, and predicted by drug-use impulsivity and estimated intelligence. These findings highlight a novel and sensitive measure of acute DD that is easy to administer.'

# Update view state: scroll so row 23 is at top, select D35 as the active cell
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
